$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new date/value pair
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 3.083829270092098

# Remove rows 3 through 17 (old trailing data no longer needed)
$ws.Range("A3:B17").EntireRow.Delete()
